$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-514).
# All of them are being bumped forward by exactly one day (46061 -> 46062).
for ($row = 2; $row -le 514; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current + 1
    }
}
